$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: replace the SUM formula row with a new component entry (LCD I2C 16x2)
$ws.Range("D9").ClearContents()
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "LCD I2C 16x2"
$ws.Cells.Item(9, 3).Value = 1

# Row 10: new component entry (Button)
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Button"
$ws.Cells.Item(10, 3).Value = 3

# Row 11: new component entry (USB UART CP2102)
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "USB UART CP2102"
$ws.Cells.Item(11, 3).Value = 1

# Apply the same style as the other data cells (style index 2 in styles.xml:
# Times New Roman 13, centered) to the newly populated cells in A9:C11 and D9
$srcStyleRange = $ws.Range("A8:C8")
$srcStyleRange.Copy()
$ws.Range("A9:C11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selected cell to D10
$ws.Range("D10").Select()
